$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86 (Cell H1): was Ladder/Ladder/0.156 -> now Sample 3920_Final, Type Sample, Dilution 80, Replicate 1
$ws.Range("B86").Value = "3920_Final"
$ws.Range("C86").Value = "Sample"
$ws.Range("D86").Value = ""
$ws.Range("E86").Value = 80
$ws.Range("F86").Value = 1

# Row 87 (Cell H2): was Ladder/Ladder/0.156 -> now Sample 3932_Final, Type Sample, Dilution 80, Replicate 1
$ws.Range("B87").Value = "3932_Final"
$ws.Range("C87").Value = "Sample"
$ws.Range("D87").Value = ""
$ws.Range("E87").Value = 80
$ws.Range("F87").Value = 1

# Row 88 (Cell H3): was 4 D6/Sample/Dilution 20/Replicate 1 -> now Sample 3936_Final, Dilution 80, Replicate 1
$ws.Range("B88").Value = "3936_Final"
$ws.Range("C88").Value = "Sample"
$ws.Range("D88").Value = ""
$ws.Range("E88").Value = 80
$ws.Range("F88").Value = 1

# Row 89 (Cell H4): was 4 D6/Sample/Dilution 20/Replicate 2 -> now all blank (B:F)
$ws.Range("B89").Value = ""
$ws.Range("C89").Value = ""
$ws.Range("D89").Value = ""
$ws.Range("E89").Value = ""
$ws.Range("F89").Value = ""

# Row 90 (Cell H5): was 8 D6/Sample/Dilution 20/Replicate 1 -> now all blank (B:F)
$ws.Range("B90").Value = ""
$ws.Range("C90").Value = ""
$ws.Range("D90").Value = ""
$ws.Range("E90").Value = ""
$ws.Range("F90").Value = ""

# Row 91 (Cell H6): was 8 D6/Sample/Dilution 20/Replicate 2 -> now all blank (B:F)
$ws.Range("B91").Value = ""
$ws.Range("C91").Value = ""
$ws.Range("D91").Value = ""
$ws.Range("E91").Value = ""
$ws.Range("F91").Value = ""

# Update the selection/view to match the final state (best-effort; engine may not persist topLeftCell)
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("E86:E88").Select()
